$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last checked" timestamp (column D) for all existing data rows (2-54)
$ws.Range("D2:D54").Value = 45999.322824074072

# New data rows (18-54): station name, terminal name, last-charge-end time
$newRows = @(
    ,@("长沙特来电飞狐四方坪南区充电站","101号直流",45987.55260416667)
    ,@("长沙特来电飞狐四方坪南区充电站","201号直流",45994.55159722222)
    ,@("长沙特来电飞狐四方坪东区充电站","005A号直流",45995.09239583334)
    ,@("长沙特来电飞狐四方坪西区充电站","503号直流",45996.55458333333)
    ,@("长沙特来电飞狐四方坪西区充电站","702号直流",45997.07376157407)
    ,@("长沙特来电飞狐四方坪西区充电站","603号直流",45997.24623842593)
    ,@("长沙特来电飞狐四方坪西区充电站","904号直流",45997.30842592593)
    ,@("长沙特来电飞狐四方坪东区充电站","003B号直流",45997.55152777778)
    ,@("长沙特来电飞狐四方坪东区充电站","905号直流",45997.55268518518)
    ,@("长沙特来电飞狐四方坪东区充电站","004A号直流",45998.04613425926)
    ,@("长沙特来电飞狐四方坪西区充电站","705号直流",45998.04956018519)
    ,@("长沙特来电飞狐四方坪西区充电站","604号直流",45998.05498842592)
    ,@("长沙特来电飞狐四方坪南区充电站","202号直流",45998.14129629629)
    ,@("长沙特来电飞狐四方坪东区充电站","002B号直流",45998.237280092595)
    ,@("长沙特来电飞狐四方坪西区充电站","405号直流",45998.23787037037)
    ,@("长沙特来电飞狐四方坪西区充电站","402号直流",45998.26613425926)
    ,@("长沙特来电飞狐四方坪南区充电站","103号直流",45998.326527777775)
    ,@("长沙特来电飞狐四方坪西区充电站","804号直流",45998.41106481481)
    ,@("长沙特来电飞狐四方坪东区充电站","402号直流",45998.54318287037)
    ,@("长沙特来电飞狐四方坪南区充电站","401号直流",45998.545949074076)
    ,@("长沙特来电飞狐四方坪南区充电站","306号直流",45998.55434027778)
    ,@("长沙特来电飞狐四方坪东区充电站","001B号直流",45998.57047453704)
    ,@("长沙特来电飞狐四方坪西区充电站","902号直流",45998.57313657407)
    ,@("长沙特来电飞狐四方坪西区充电站","803号直流",45998.5808912037)
    ,@("长沙特来电飞狐四方坪西区充电站","703号直流",45998.580983796295)
    ,@("长沙特来电飞狐四方坪西区充电站","505号直流",45998.584502314814)
    ,@("长沙特来电飞狐四方坪西区充电站","B01号直流",45998.59831018518)
    ,@("长沙特来电飞狐四方坪东区充电站","006B号直流",45998.63952546296)
    ,@("长沙特来电飞狐四方坪南区充电站","105号直流",45998.65846064815)
    ,@("长沙特来电飞狐四方坪西区充电站","802号直流",45998.67444444444)
    ,@("长沙特来电飞狐四方坪西区充电站","805号直流",45998.677569444444)
    ,@("长沙特来电飞狐四方坪西区充电站","905号直流",45998.689479166664)
    ,@("长沙市开福区高岭香江国际城充电站建设项目","108号直流",45998.70358796296)
    ,@("长沙特来电飞狐四方坪东区充电站","401号直流",45998.71167824074)
    ,@("长沙市开福区高岭香江国际城充电站建设项目","309号直流",45998.721041666664)
    ,@("长沙特来电飞狐四方坪西区充电站","903号直流",45998.7409837963)
    ,@("长沙特来电飞狐四方坪东区充电站","009A号直流",45998.789293981485)
)

for ($i = 0; $i -lt $newRows.Length; $i++) {
    $r = 18 + $i
    $row = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
}

# Update selection to reflect the cell active when the workbook was saved
$ws.Range("E23").Select()
